$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 cells to the new evaluation run values
$ws.Range("A2").Value = 45656.97892751805
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.8
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 1
$ws.Range("V2").Value = 2.8
$ws.Range("AA2").Value = 1
$ws.Range("AD2").Value = 2.8
$ws.Range("AI2").Value = 3
$ws.Range("AL2").Value = 5.600000000000001
$ws.Range("AQ2").Value = 8.600000000000001
$ws.Range("AS2").Value = 8.600000000000001
$ws.Range("AT2").Value = "The opening question was basic, expecting a polite response without eliciting direct job-to-be-done information."
$ws.Range("AU2").Value = "The persona maintained a pleasant demeanor helping to create rapport.; The response was professionally tailored to his work context."
$ws.Range("AV2").Value = "The initial question could be more specific or purposive to inquire about the personas's work or experience."
$ws.Range("AW2").Value = "Can you tell me about a recent challenge in digital banking?; What aspects of risk management are you currently focusing on?"

# Remove the old row 3 (second interview Q&A entry) entirely
$ws.Rows.Item(3).Delete()

